$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Success"
